$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.734.87'
$ws.Range('E2').Value = '  +0.22%  '

$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').Value = '211.63'
$ws.Range('E5').Value = '  -0.06%  '

$ws.Range('E6').Value = '  -0.48%  '

$ws.Range('E7').Value = '  +0.21%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').Value = '0.248'
$ws.Range('E9').Value = '  +0.45%  '

$ws.Range('D10').Value = '19.77'

$ws.Range('E11').Value = '  +0.69%  '

$ws.Range('D12').Value = '1.826.83'
$ws.Range('E12').Value = '  +0.16%  '

$ws.Range('D13').Value = '1.599.83'
$ws.Range('E13').Value = '  +0.62%  '

$ws.Range('E14').Value = '  +0.36%  '

$ws.Range('D15').Value = '0.522'
$ws.Range('E15').Value = '  -0.31%  '

$ws.Range('D16').Value = '65.06'
$ws.Range('E16').Value = '  -0.07%  '

$ws.Range('D17').Value = '0.0₃0740'
$ws.Range('E17').Value = '  +0.35%  '

$ws.Range('D18').Value = '210.43'
$ws.Range('E18').Value = '  +0.85%  '

$ws.Range('E19').Value = '  +0.20%  '

$ws.Range('E20').Value = '  +1.99%  '

$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('E22').Value = '  -2.02%  '

$ws.Range('D23').Value = '8.99'
$ws.Range('E23').Value = '  +0.12%  '

$ws.Range('D24').Value = '143.69'
$ws.Range('E24').Value = '  -1.06%  '

$ws.Range('E25').Value = '  +0.14%  '

$ws.Range('E26').Value = '  -0.28%  '

$ws.Range('E27').Value = '  -0.91%  '

$ws.Range('E28').Value = '  +0.56%  '

$ws.Range('E29').Value = '  -0.56%  '

$ws.Range('E30').Value = '  -0.22%  '

$ws.Range('E31').Value = '  +1.20%  '

$ws.Range('E32').Value = '  +1.20%  '

$ws.Range('D33').Value = '1.294.72'
$ws.Range('E33').Value = '  +1.40%  '

$ws.Range('E34').Value = '  +0.82%  '

$ws.Range('E35').Value = '  +0.68%  '

$ws.Range('D36').Value = '0.604'
$ws.Range('E36').Value = '  -2.68%  '

$ws.Range('D37').Value = '1.15'
$ws.Range('E37').Value = '  +10.70%  '

$ws.Range('D38').Value = '0.0170'
$ws.Range('E38').Value = '  -0.30%  '

$ws.Range('D39').Value = '0.831'
$ws.Range('E39').Value = '  -0.49%  '

$ws.Range('E40').Value = '  -1.79%  '

$ws.Range('E41').Value = '  -0.16%  '

$ws.Range('D42').Value = '0.785'
$ws.Range('E42').Value = '  +0.12%  '

$ws.Range('D43').Value = '63.04'
$ws.Range('E43').Value = '  -1.55%  '

$ws.Range('D44').Value = '1.738.89'
$ws.Range('E44').Value = '  +0.23%  '

$ws.Range('D45').Value = '90.70'
$ws.Range('E45').Value = '  -0.57%  '

$ws.Range('E46').Value = '  -2.16%  '

$ws.Range('E47').Value = '  -0.48%  '

$ws.Range('D48').Value = '0.0518'
$ws.Range('E48').Value = '  +1.94%  '

$ws.Range('E49').Value = '  +0.14%  '

$ws.Range('D50').Value = '7.39'
$ws.Range('E50').Value = '  +0.17%  '

$ws.Range('D51').Value = '0.396'
$ws.Range('E51').Value = '  +0.88%  '
